$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("06-09-2021", 3.55, 3.19, 3.21),
    @("07-09-2021", 3.59, 3.2, 3.23),
    @("08-09-2021", 3.47, 3.26, 3.27),
    @("09-09-2021", 3.44, 3.33, 3.3),
    @("10-09-2021", 3.42, 3.38, 3.3)
)

$startRow = 175
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $cellA = $ws.Cells.Item($row, 1)
    # Force the date-looking string to be stored as literal text (shared
    # string) instead of letting Excel auto-convert it to a date serial.
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
